# report: updated values in table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Updated statistics values in the three result tables ----

# Table 1 (rows 7-11)
$ws.Range("D8").Value = 153
$ws.Range("E8").Value = 0.79
$ws.Range("F8").Value = 30
$ws.Range("G8").Value = 20
$ws.Range("H8").Value = 25.132999999999999

$ws.Range("D10").Value = 145
$ws.Range("E10").Value = 0.74
$ws.Range("G10").Value = 24
$ws.Range("H10").Value = 28.867000000000001

$ws.Range("D11").Value = 153
$ws.Range("E11").Value = 0.79
$ws.Range("G11").Value = 20
$ws.Range("H11").Value = 25.132999999999999

# Table 2 (rows 17-21)
$ws.Range("D18").Value = 141
$ws.Range("E18").Value = 0.59
$ws.Range("F18").Value = 45
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = 33.732999999999997

$ws.Range("D20").Value = 196
$ws.Range("E20").Value = 0.72
$ws.Range("F20").Value = 45
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 36.567

$ws.Range("D21").Value = 141
$ws.Range("G21").Value = 29
$ws.Range("H21").Value = 33.732999999999997

# Table 3 (rows 27-31)
$ws.Range("D28").Value = 135
$ws.Range("E28").Value = 0.39
$ws.Range("F28").Value = 51
$ws.Range("G28").Value = 25
$ws.Range("H28").Value = 29.533000000000001

$ws.Range("D31").Value = 151
$ws.Range("E31").Value = 0.45
$ws.Range("F31").Value = 51
$ws.Range("G31").Value = 25
$ws.Range("H31").Value = 30.067

# ---- H11 / H20 no longer carry the highlighted-fill format; match the
#      plain style already used by their neighboring rows (e.g. H10/H18) ----
$ws.Range("H10").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H18").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Selection / view state ----
$ws.Range("D33:D34").Select()
$excel.ActiveWindow.Zoom = 125
